# Generate Report for Handoff
# Adds a new handed-off file (ae96f33b-d1d8-457a-af9b-f05a5433f551) as row 3
# on the Overview / zh-cn / de-de sheets, mirroring the existing
# 78e2c296-f1dc-4bd8-95b2-65a4a58df23b row.

$wb = $excel.ActiveWorkbook

$fileName   = "ae96f33b-d1d8-457a-af9b-f05a5433f551.md"
$zhXlf      = "ae96f33b-d1d8-457a-af9b-f05a5433f551.1c6c5cca8ee3de439c45464cf54f5be0b75b0dca.zh-cn.xlf"
$deXlf      = "ae96f33b-d1d8-457a-af9b-f05a5433f551.1c6c5cca8ee3de439c45464cf54f5be0b75b0dca.de-de.xlf"
$status     = "Ready for handoff"
$dtOverview = "2016-03-22 14:38:54"
$dtZh       = "2016-03-22 14:38:50"
$dtDe       = "2016-03-22 14:38:54"
$epoch      = "0001-01-01 00:00:00"
$ext        = ".md"
$reason     = "Include"
$dateFmt    = "yyyy-mm-dd HH:mm:ss"
$linkColor  = 15570276  # OLE (BGR) form of RGB 6495ED, matches existing hyperlink font color

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/3eb9bab449976525ce0b8e768563ef21aaa9743d/e2e/ae96f33b-d1d8-457a-af9b-f05a5433f551.md"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50afbfa611627b0b8ce703319bf21f0dd358d6da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ae96f33b-d1d8-457a-af9b-f05a5433f551.1c6c5cca8ee3de439c45464cf54f5be0b75b0dca.zh-cn.xlf"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d2bb3f35981fe799c8ce4bac8ec5cd24fe35985/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ae96f33b-d1d8-457a-af9b-f05a5433f551.1c6c5cca8ee3de439c45464cf54f5be0b75b0dca.de-de.xlf"

function Style-LinkCell($cell) {
    $cell.Font.Underline = 2
    $cell.Font.Color = $linkColor
}

function Style-DateCell($cell) {
    $cell.NumberFormat = $dateFmt
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1): File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = $dtOverview

Style-LinkCell($wsOverview.Range("A3"))
Style-DateCell($wsOverview.Range("D3"))

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl, "", "", $fileName)

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $fileName
$wsZh.Range("B3").Value = $ext
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = $zhXlf
$wsZh.Range("E3").Value = $dtZh
$wsZh.Range("H3").Value = $epoch
$wsZh.Range("J3").Value = $reason

Style-LinkCell($wsZh.Range("A3"))
Style-LinkCell($wsZh.Range("D3"))
Style-DateCell($wsZh.Range("E3"))
Style-DateCell($wsZh.Range("H3"))

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $fileName)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhUrl, "", "", $zhXlf)

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $fileName
$wsDe.Range("B3").Value = $ext
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = $deXlf
$wsDe.Range("E3").Value = $dtDe
$wsDe.Range("H3").Value = $epoch
$wsDe.Range("J3").Value = $reason

Style-LinkCell($wsDe.Range("A3"))
Style-LinkCell($wsDe.Range("D3"))
Style-DateCell($wsDe.Range("E3"))
Style-DateCell($wsDe.Range("H3"))

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $fileName)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deUrl, "", "", $deXlf)
